# Apply updated crypto price/volume figures (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value still parses as a plain number need to be
# pinned to Text format first, otherwise COM auto-coerces them to a Double
# and the decimal-grouped "price" strings (e.g. "241.61") lose their
# original text representation. Style is reset back to Normal afterwards so
# no visible formatting change is introduced.
$textCells = @('D5', 'D6', 'D8', 'D9', 'D10', 'D11', 'D13', 'D14', 'D15', 'D17', 'D18', 'D19', 'D23', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D44', 'D47', 'D49', 'D50', 'D51')
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = '@'
}

$ws.Range('D2').Value = '29.285.59'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').Value = '1.871.83'
$ws.Range('E3').Value = '  +0.16%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '0.7091'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').Value = '241.61'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '0.07816'
$ws.Range('E8').Value = '  +1.94%  '
$ws.Range('D9').Value = '0.3096'
$ws.Range('E9').Value = '  -0.71%  '
$ws.Range('D10').Value = '25.00'
$ws.Range('E10').Value = '  +1.21%  '
$ws.Range('D11').Value = '0.08412'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '1.875.44'
$ws.Range('E12').Value = '  +0.04%  '
$ws.Range('D13').Value = '5.231'
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').Value = '0.7107'
$ws.Range('E14').Value = '  -0.02%  '
$ws.Range('D15').Value = '91.00'
$ws.Range('E15').Value = '  -0.30%  '
$ws.Range('D16').Value = '29.292.01'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '6.072'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '0.000008189'
$ws.Range('E18').Value = '  +4.09%  '
$ws.Range('D19').Value = '239.84'
$ws.Range('E19').Value = '  -1.58%  '
$ws.Range('D21').Value = '2.123.16'
$ws.Range('E21').Value = '  +0.45%  '
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').Value = '7.734'
$ws.Range('E23').Value = '  -1.38%  '
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').Value = '0.1595'
$ws.Range('E25').Value = '  -2.40%  '
$ws.Range('D26').Value = '162.60'
$ws.Range('E26').Value = '  -0.68%  '
$ws.Range('D27').Value = '9.008'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('E28').Value = '  -0.34%  '
$ws.Range('D29').Value = '1.503'
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('D30').Value = '4.392'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('D31').Value = '1.295'
$ws.Range('E31').Value = '  -0.76%  '
$ws.Range('D32').Value = '4.301'
$ws.Range('E32').Value = '  +1.45%  '
$ws.Range('D33').Value = '0.05357'
$ws.Range('E33').Value = '  +3.81%  '
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('D36').Value = '0.7483'
$ws.Range('E36').Value = '  -4.24%  '
$ws.Range('D37').Value = '2.692'
$ws.Range('E37').Value = '  +0.29%  '
$ws.Range('D38').Value = '0.01871'
$ws.Range('E38').Value = '  +0.90%  '
$ws.Range('D39').Value = '1.219.76'
$ws.Range('E39').Value = '  +5.39%  '
$ws.Range('D40').Value = '2.725'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').Value = '6.492'
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('D42').Value = '0.8888'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').Value = '108.67'
$ws.Range('E44').Value = '  +5.09%  '
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '2.021.19'
$ws.Range('E46').Value = '  +0.39%  '
$ws.Range('D47').Value = '0.5201'
$ws.Range('E48').Value = '  +2.13%  '
$ws.Range('D49').Value = '1.795'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').Value = '9.396'
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').Value = '0.4311'
$ws.Range('E51').Value = '  +0.36%  '

foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = 'Normal'
}
